# Merge the split "<id>...</id>" runs (three separate runs: "<id>",
# the bare id text, "</id>") into a single run per occurrence, so the
# whole tag ends up as one run using the <id>/</id> run's formatting
# (Courier New, color 7f6000). This mirrors newly-downloaded tc/tcn/tl
# content where the id tag is written as a single literal run.

$d = $word.ActiveDocument

$ids = @("p157v_1", "p157v_2", "p157v_3")

foreach ($id in $ids) {
    $needle = "<id>" + $id + "</id>"
    $d.Content.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, $needle, 2) | Out-Null
}
